$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.01979999999999
$ws.Range("A7").Value = -19.78149999999999
$ws.Range("C7").Value = -12.58270000000001
$ws.Range("C15").Value = -14.76569999999999
$ws.Range("A16").Value = -21.64549999999999
$ws.Range("D16").Value = -9.054900000000004
$ws.Range("D19").Value = -8.509099999999995
$ws.Range("C21").Value = -12.30950000000001
$ws.Range("C22").Value = -12.51760000000001
$ws.Range("C23").Value = -12.3891
$ws.Range("A28").Value = -22.0264
$ws.Range("A29").Value = -21.03409999999997
$ws.Range("A32").Value = -21.20219999999998
$ws.Range("C34").Value = -11.88640000000001
$ws.Range("D36").Value = -8.513299999999994
$ws.Range("A40").Value = -20.34619999999999
$ws.Range("C43").Value = -13.2596
$ws.Range("C45").Value = -13.9682
$ws.Range("D46").Value = -8.274599999999996
$ws.Range("C50").Value = -14.36879999999999
$ws.Range("D50").Value = -8.02
$ws.Range("C51").Value = -12.18149999999999
$ws.Range("A52").Value = -22.2165
$ws.Range("A57").Value = -22.39280000000002
$ws.Range("A66").Value = -22.029
$ws.Range("C66").Value = -12.4051
$ws.Range("C67").Value = -11.0741
$ws.Range("C79").Value = -11.55850000000001
$ws.Range("C84").Value = -12.95479999999999
$ws.Range("C92").Value = -11.39370000000001
$ws.Range("D95").Value = -8.333899999999998
$ws.Range("C97").Value = -11.29860000000001
$ws.Range("D97").Value = -8.416899999999996
$ws.Range("A100").Value = -21.9704
